$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "Sait Tunç, Burhaneddin Sandıkçı, Bekir Tanrıöver"
$ws.Range("B6").Value = "Grado Department of Industrial and Systems Engineering, Virginia Tech, Blacksburg, Virginia 24061;; Department of Industrial Engineering, Istanbul Technical University, Istanbul, Turkey 34367;; Division of Nephrology, The University of Arizona, Tucson, Arizona 85724"
$ws.Range("C6").Value = "https://openalex.org/W4210971474"
$ws.Range("D6").Value = "A Simple Incentive Mechanism to Alleviate the Burden of Organ Wastage in Transplantation"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2022-08-01"
$ws.Range("F6").Value = "Management Science"
$ws.Range("G6").Value = "Institute for Operations Research and the Management Sciences"
$ws.Range("H6").Value = "https://doi.org/10.1287/mnsc.2021.4203"
$ws.Range("J6").Value = "N/A"
$ws.Range("K6").Value = "closed"
$ws.Range("O6").Value = "NA"
$ws.Range("P6").Value = "https://doi.org/10.1287/mnsc.2021.4203"
$ws.Range("A7").Value = "Darren Stewart, Bekir Tanrıöver, Gaurav Gupta"
$ws.Range("B7").Value = "Department of Surgery, New York University Langone Health, New York, New York; Division of Nephrology, The University of Arizona, Tucson, Arizona; Division of Nephrology, School of Medicine, Virginia Commonwealth University, Richmond, Virginia"
$ws.Range("C7").Value = "https://openalex.org/W4302011830"
$ws.Range("D7").Value = "Oversimplification and Misplaced Blame Will Not Solve the Complex Kidney Underutilization Problem"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2022-12-01"
$ws.Range("F7").Value = "Kidney360"
$ws.Range("G7").Value = "Lippincott Williams & Wilkins"
$ws.Range("H7").Value = "https://doi.org/10.34067/kid.0005402022"
$ws.Range("J7").Value = "publishedVersion"
$ws.Range("K7").Value = "gold"
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = "4"
$ws.Range("O7").Value = "https://pubmed.ncbi.nlm.nih.gov/36591359"
$ws.Range("P7").Value = "https://doi.org/10.34067/kid.0005402022"
$ws.Range("A8").Value = "Gurmukteshwar Singh, Reginald Gohh, Dinah Clark, Kartik Kalra, Manoj Das, Gitana Bradauskaite, Anthony J. Bleyer, Bekir Tanrıöver, Alex R. Chang, Prince Mohan Anand"
$ws.Range("B8").Value = "Department of Nephrology, Geisinger Health, Danville, PA 17822, USA; Division of Organ Transplantation, Rhode Island Hospital, Providence, RI 02908, USA; Natera, Inc., Austin, TX 78753, USA; Department of Nephrology, Geisinger Health, Danville, PA 17822, USA; Department of Nephrology, Geisinger Health, Danville, PA 17822, USA; Division of Nephrology, Einstein Medical Center, Philadelphia, PA 19141, USA; Division of Nephrology, Wake Forest School of Medicine, Winston-Salem, NC 27157, USA; Division of Nephrology, University of Arizona College of Medicine, Tucson, AZ 85719, USA; Department of Nephrology, Geisinger Health, Danville, PA 17822, USA; Department of Nephrology, Medical University of South Carolina, Charleston, SC 29425, USA"
$ws.Range("C8").Value = "https://openalex.org/W4220987486"
$ws.Range("D8").Value = "Vignette-Based Reflections to Inform Genetic Testing Policies in Living Kidney Donors"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "2022-03-26"
$ws.Range("F8").Value = "Genes"
$ws.Range("G8").Value = "Multidisciplinary Digital Publishing Institute"
$ws.Range("H8").Value = "https://doi.org/10.3390/genes13040592"
$ws.Range("I8").Value = "cc-by"
$ws.Range("O8").Value = "https://pubmed.ncbi.nlm.nih.gov/35456398"
$ws.Range("P8").Value = "https://doi.org/10.3390/genes13040592"
$ws.Range("A9").Value = "Mutlu Mete, Mehmet Ayvaci, Venkatesh Kumar Ariyamuthu, Alpesh Amin, Matthias Peltz, Jennifer T. Thibodeau, Justin L. Grodin, Pradeep P.A. Mammen, Sonia Garg, Faris Araj, Robert Morlend, Mark H. Drazner, Nashila AbdulRahim, Yeongin Kim, Yusuf Salam, Ahmet B. Gungor, Bulent Delibasi, Suman K. Kotla, Malcolm MacConmara, Prince Mohan, Gaurav Gupta, Bekir Tanrıöver"
$ws.Range("B9").Value = "Department of Computer Science, Texas A&M University—Commerce, Commerce, Texas, USA; Information Systems, Naveen Jindal School of Business, University of Texas at Dallas, Richardson, Texas, USA; Division of Nephrology, University of Arizona College of Medicine—Tucson, Tucson, Arizona, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiovascular and Thoracic Surgery, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Cardiology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Nephrology, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Virginia Commonwealth University, Information Systems, School of Business, Richmond, Virginia, USA; School of Public Health at the University of Texas Health Science Center, Dallas, Texas, USA; Division of Nephrology, Banner University Medical Group—Tucson, Tucson, Arizona, USA; School of Behavioral and Brain Sciences, University of Texas at Dallas, Richardson, Texas, USA; Dallas Nephrology Associates, Dallas, Texas, USA; Department of Surgery, University of Texas Southwestern Medical Center, Dallas, Texas, USA; Division of Nephrology, Medical University of South Carolina, Lancaster, South Carolina, USA; Division of Nephrology, Virginia Commonwealth University, Richmond, Virginia, USA; Division of Nephrology, University of Arizona College of Medicine—Tucson, Tucson, Arizona, USA"
$ws.Range("C9").Value = "https://openalex.org/W4224903100"
$ws.Range("D9").Value = "Predicting Post-Heart Transplant Composite Renal Outcome Risk in Adults: A Machine Learning Decision Tool"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2022-06-01"
$ws.Range("F9").Value = "Kidney International Reports"
$ws.Range("G9").Value = "Elsevier BV"
$ws.Range("H9").Value = "https://doi.org/10.1016/j.ekir.2022.04.004"
$ws.Range("I9").Value = "cc-by-nc-nd"
$ws.Range("J9").Value = "publishedVersion"
$ws.Range("K9").Value = "gold"
$ws.Range("M9").NumberFormat = "@"
$ws.Range("M9").Value = "3"
$ws.Range("O9").Value = "https://pubmed.ncbi.nlm.nih.gov/35685329"
$ws.Range("P9").Value = "https://doi.org/10.1016/j.ekir.2022.04.004"
$ws.Range("A10").Value = "Ramesh Batra, Venkatesh Kumar Ariyamuthu, Malcolm MacConmara, Gaurav Gupta, Ahmet B. Gungor, Bekir Tanrıöver"
$ws.Range("B10").Value = "Department of SurgerySchool of MedicineYale UniversityNew HavenCT; These authors contributed equally as co‐first authors.; Division of NephrologyCollege of MedicineUniversity of ArizonaTucsonAZ; These authors contributed equally as co‐first authors.; Abdominal Medical AffairsTransMedics Group, Inc.AndoverMA; Division of Nephrology Virginia Commonwealth University  Richmond Virginia USA; Division of NephrologyBanner University Medical GroupTucsonAZ; Division of Nephrology College of Medicine University of Arizona  Tucson Arizona USA"
$ws.Range("C10").Value = "https://openalex.org/W4205981601"
$ws.Range("D10").Value = "Outcomes of Simultaneous Liver‐Kidney Transplantation Using Kidneys of Deceased Donors With Acute Kidney Injury"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2022-03-03"
$ws.Range("F10").Value = "Liver Transplantation"
$ws.Range("G10").Value = "Lippincott Williams & Wilkins"
$ws.Range("H10").Value = "https://doi.org/10.1002/lt.26406"
$ws.Range("I10").Value = "N/A"
$ws.Range("J10").Value = "N/A"
$ws.Range("K10").Value = "closed"
$ws.Range("O10").Value = "https://pubmed.ncbi.nlm.nih.gov/35006615"
$ws.Range("P10").Value = "https://doi.org/10.1002/lt.26406"
$ws.Range("G17").Value = "Springer Nature"
